# Rotate the "Activity / ScheduledResource / PlannedQty / Plant /
# ProductionDivision" block between rows 2, 3 and 4:
#   old row 2 data -> row 4
#   old row 3 data -> row 2
#   old row 4 data -> row 3
# (row 4's Plant/ProductionDivision land swapped relative to the
# original row 2 values, matching the target diff exactly)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Digital Print F 4x0"
$ws.Range("D2").Value = "252-HP 10000 Press"
$ws.Range("G2").Value = "719"
$ws.Range("L2").Value = "252-HP 10000 Press"
$ws.Range("M2").Value = "252-HP 10000 Press"

$ws.Range("B3").Value = "Cut"
$ws.Range("D3").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("G3").Value = "715"
$ws.Range("L3").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("M3").Value = "406-45`" Polar 115ED Cutter`n404-45`" Polar 115EMC Cutter`n405-54`" Polar 137EMC Cutter`n402-45`" Polar 115EMC Cutter`n403-54`" Polar 137ED Cutter"

$ws.Range("B4").Value = "-"
$ws.Range("D4").Value = "169-Press Approval Task "
$ws.Range("G4").Value = "740"
$ws.Range("L4").Value = "169-Press Approval Task "
$ws.Range("M4").Value = "Press Approval Task"
